$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input values in row 4
$ws.Range("D4").Value = 1560
$ws.Range("F4").Value = 2000
$ws.Range("G4").Value = 2120

# Update input values in row 6
$ws.Range("D6").Value = 1720
$ws.Range("G6").Value = 2150

# Update the selected cell on the sheet
$ws.Range("K13").Select()
